# Auto-generated Excel COM-interop script
# Applies a scheduled market-price data refresh to the per-job Golem profit sheets
# (columns H:N -- currentAveragePrice*, LevePrice*, LeveProfit*) across all 8 job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 70.44444
$ws.Range("I6").Value = 31.75
$ws.Range("J6").Value = 380
$ws.Range("K6").Value = 95.25
$ws.Range("L6").Value = 1140
$ws.Range("M6").Value = 16.75
$ws.Range("N6").Value = -1364
# Row 11
$ws.Range("H11").Value = 275.07144
$ws.Range("I11").Value = 275.07144
$ws.Range("K11").Value = 275.07144
$ws.Range("M11").Value = -135.07144
# Row 12
$ws.Range("H12").Value = 1099.6
# Row 28
$ws.Range("H28").Value = 1736.8462
$ws.Range("I28").Value = 617.8570999999999
$ws.Range("J28").Value = 3042.3333
$ws.Range("K28").Value = 617.8570999999999
$ws.Range("L28").Value = 3042.3333
$ws.Range("M28").Value = -132.8570999999999
$ws.Range("N28").Value = -4012.3333
# Row 33
$ws.Range("H33").Value = 407.24
$ws.Range("I33").Value = 382.54166
$ws.Range("K33").Value = 382.54166
$ws.Range("M33").Value = -153.54166
# Row 38
$ws.Range("H38").Value = 256.75
$ws.Range("J38").Value = 882
$ws.Range("L38").Value = 2646
$ws.Range("N38").Value = -3390
# Row 39
$ws.Range("H39").Value = 2539.75
$ws.Range("I39").Value = 1047.8
$ws.Range("K39").Value = 3143.4
$ws.Range("M39").Value = -2847.4
# Row 98
$ws.Range("H98").Value = 1506.7693
$ws.Range("I98").Value = 1506.7693
$ws.Range("K98").Value = 1506.7693
$ws.Range("M98").Value = -8.76929999999993
# Row 107
$ws.Range("H107").Value = 57051.375
$ws.Range("I107").Value = 57051.375
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 57051.375
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -55131.375
$ws.Range("N107").ClearContents()
# Row 122
$ws.Range("H122").Value = 1506.7693
$ws.Range("I122").Value = 1506.7693
$ws.Range("K122").Value = 4520.3079
$ws.Range("M122").Value = -2070.3079
# Row 138
$ws.Range("H138").Value = 3210.261
$ws.Range("I138").Value = 931
$ws.Range("J138").Value = 4014.7058
$ws.Range("K138").Value = 2793
$ws.Range("L138").Value = 12044.1174
$ws.Range("M138").Value = 2347
$ws.Range("N138").Value = -22324.1174

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 374.8
$ws.Range("I2").Value = 362.3
$ws.Range("K2").Value = 362.3
$ws.Range("M2").Value = -249.3
# Row 74
$ws.Range("H74").Value = 3925
$ws.Range("I74").Value = 2950
$ws.Range("K74").Value = 2950
$ws.Range("M74").Value = -2076
# Row 77
$ws.Range("H77").Value = 3925
$ws.Range("I77").Value = 2950
$ws.Range("K77").Value = 14750
$ws.Range("M77").Value = -10382
# Row 116
$ws.Range("H116").Value = 374.8
$ws.Range("I116").Value = 362.3
$ws.Range("K116").Value = 362.3
$ws.Range("M116").Value = 1931.7
# Row 122
$ws.Range("H122").Value = 9500
$ws.Range("I122").Value = 9500
$ws.Range("K122").Value = 28500
$ws.Range("M122").Value = -26050
# Row 125
$ws.Range("H125").Value = 82500
$ws.Range("J125").Value = 82500
$ws.Range("L125").Value = 82500
$ws.Range("N125").Value = -92340

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 374.8
$ws.Range("I3").Value = 362.3
$ws.Range("K3").Value = 362.3
$ws.Range("M3").Value = -248.3
# Row 107
$ws.Range("H107").Value = 31497.076
$ws.Range("J107").Value = 1167
$ws.Range("L107").Value = 1167
$ws.Range("N107").Value = -5007

$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 5000
$ws.Range("K55").Value = 5000
$ws.Range("M55").Value = -4685
# Row 88
$ws.Range("H88").Value = 43868.6
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 43868.6
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 43868.6
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -44680.6
# Row 91
$ws.Range("H91").Value = 43868.6
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 43868.6
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 43868.6
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -46676.6
# Row 107
$ws.Range("H107").Value = 671.53845
$ws.Range("I107").Value = 622
$ws.Range("K107").Value = 622
$ws.Range("M107").Value = 1298
# Row 132
$ws.Range("H132").Value = 993.3333
$ws.Range("I132").Value = 993.3333
$ws.Range("K132").Value = 2979.9999
$ws.Range("M132").Value = -449.9998999999998
# Row 134
$ws.Range("H134").Value = 3131
$ws.Range("J134").Value = 9927.5
$ws.Range("L134").Value = 29782.5
$ws.Range("N134").Value = -34852.5

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 794.86487
$ws.Range("I44").Value = 168.5
$ws.Range("K44").Value = 505.5
$ws.Range("M44").Value = -107.5
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
# Row 99
$ws.Range("H99").Value = 7500
$ws.Range("I99").Value = 7500
$ws.Range("K99").Value = 22500
$ws.Range("M99").Value = -20254
# Row 107
$ws.Range("H107").Value = 1301.5
$ws.Range("J107").Value = 1401
$ws.Range("L107").Value = 4203
$ws.Range("N107").Value = -8043

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 41667884
$ws.Range("I107").Value = 548.6
$ws.Range("K107").Value = 548.6
$ws.Range("M107").Value = 1371.4
# Row 136
$ws.Range("H136").Value = 19785.5
$ws.Range("J136").Value = 19785.5
$ws.Range("L136").Value = 59356.5
$ws.Range("N136").Value = -64456.5

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4962.625
$ws.Range("I122").Value = 3559.4
$ws.Range("K122").Value = 10678.2
$ws.Range("M122").Value = -8228.200000000001
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 569.2
$ws.Range("I113").Value = 442
$ws.Range("K113").Value = 1326
$ws.Range("M113").Value = 844
# Row 132
$ws.Range("H132").Value = 2207
$ws.Range("I132").Value = 888.8
$ws.Range("K132").Value = 2666.4
$ws.Range("M132").Value = -136.3999999999996

